$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "264.85"
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "1.47%"

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "26.44"
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = "-2.45%"

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "4.693"
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = "-0.08%"

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "0.06091"
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = "-1.44%"

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "6.753"
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = "1.09%"

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.8519"
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = "0.08%"

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.9106"
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = "-0.79%"

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.1417"
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = "0.61%"

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.04962"
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = "6.08%"

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.07115"
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = "0.35%"

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.03145"
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = "-0.26%"

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.09029"
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = "-0.13%"

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.001543"
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = "0.49%"

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.0006064"
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = "-1.57%"

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.006024"
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = "-0.61%"

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "3.449"
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = "-0.05%"

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "3.168"
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = "0.11%"

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "2.175"
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = "-0.18%"

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "0.1282"
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = "-1.33%"

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "4.084"
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = "0.05%"

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "0.04235"
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = "-0.17%"

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "0.001183"
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = "-2.66%"

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "0.004057"
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = "6.70%"

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "0.0001201"
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = "0.09%"

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "0.0001683"
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = "6.65%"

$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = "1.18%"

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.1114"
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = "0.23%"

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.004179"
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = "2.46%"

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.002107"
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = "-3.51%"

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.01164"
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = "-28.67%"

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.00005098"
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = "-1.18%"

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.00000000751"
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = "0.11%"

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.2584"
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = "53.68%"

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.00002103"
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = "0.11%"

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.0002003"
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = "0.11%"
